# "modifisert dbinit etter den nye excel forslag"
# Add a new "maks plasser" (max seats) column (I) to the timetable, fill it
# in for the three routes that already specify a departure/arrival time
# pair, and drop a couple of now-superfluous time/duration values that were
# left over from the old layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I.
$ws.Range("I1").Value = "maks plasser"

# Sandefjord -> Strømstad: max seats for this departure.
$ws.Range("I14").Value = 500

# Strømstad -> Sandefjord: the old return-time value is no longer used.
$ws.Range("D15").ClearContents()

# Stavanger -> Bergen: max seats for this departure.
$ws.Range("I17").Value = 700

# Bergen -> Stavanger: the old return-time value is no longer used.
$ws.Range("D18").ClearContents()

# Oslo -> Kiel: max seats for this departure.
$ws.Range("I20").Value = 1000

# Kiel -> Oslo: the old duration value is no longer used.
$ws.Range("D21").ClearContents()

# Leave the selection where the author ended up editing.
$ws.Range("K14").Select()
